$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" header text in F1 (09:15 -> 09:30)
$ws.Range("F1").Value = "Last status check on: 20.01.2022 09:30"

# Row 8 (Benzina Albert Modrice): convert D8/E8 from text to real numeric values
$ws.Range("D8").Value = 0.3
$ws.Range("E8").Value = 44581.385625
$ws.Range("E8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
